$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the "Meta description: ..." paragraph that currently sits
#    right after the "Play Black Widow Slot Game for Free | Review"
#    Heading1 paragraph near the top of the document.
# ------------------------------------------------------------------
$metaPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Meta description:*") {
        $metaPara = $p
        break
    }
}
if ($metaPara -ne $null) {
    $metaPara.Range.Delete()
}

# ------------------------------------------------------------------
# 2) Insert a new paragraph just before the final paragraph (the one
#    that currently holds the "Create a feature image ..." prompt),
#    containing a bold run with the title text.
# ------------------------------------------------------------------
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphBefore()

$newParaIndex = $d.Paragraphs.Count - 1
$newPara = $d.Paragraphs.Item($newParaIndex)

$titleXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Black Widow Slot Game for Free | Review</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$newPara.Range.InsertXML($titleXml)

# ------------------------------------------------------------------
# 3) Replace the text of the final "Create a feature image ..."
#    paragraph with the meta-description wording, keeping its
#    existing italic run formatting intact. Range.Text (instead of
#    Find.Execute's replace) is used so straight apostrophes are not
#    auto-corrected into curly/smart quotes.
# ------------------------------------------------------------------
$newText = "Read our review of Black Widow slot game, available to play for free. Get a chance to win huge prizes with the game's free spins bonus round."

$finalPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$finalRange = $finalPara.Range
$finalRange.MoveEnd(1, -1)
$finalRange.Text = $newText
